# Insert a new weekly data record at row 62, shifting the existing
# rows 62-130 down to rows 63-131 (dimension grows from A1:R130 to A1:R131).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A62").EntireRow.Insert()

$ws.Range("A62").Value = 4
$ws.Range("B62").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C62").Value = "Los Lagos"
$ws.Range("D62").Value = 44789
$ws.Range("E62").Value = 10
$ws.Range("F62").Value = 100112022
$ws.Range("G62").Value = "Arveja Verde"
$ws.Range("H62").Value = "Perfection"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 70
$ws.Range("K62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("M62").Value = 45000
$ws.Range("N62").Value = "`$/malla 25 kilos"
$ws.Range("O62").Value = "Provincia de Huasco"
$ws.Range("P62").Value = 1800
$ws.Range("Q62").Value = 25
$ws.Range("R62").Value = "Hortaliza"
